$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 19 - "rom" playing "AgainstTime" (אחד נגד השני)
$ws.Cells.Item(19, 1).Value = "rom"
$ws.Cells.Item(19, 2).Value = 43989.85160238426
$ws.Cells.Item(19, 2).NumberFormat = "dd-MM-yyyy"
$ws.Cells.Item(19, 3).Value = "אחד נגד השני"
$ws.Cells.Item(19, 4).Value = 40

# Row 20 - "ליאת" playing "AgainstTime" (אחד נגד השני)
$ws.Cells.Item(20, 1).Value = "ליאת"
$ws.Cells.Item(20, 2).Value = 43989.85160335648
$ws.Cells.Item(20, 2).NumberFormat = "dd-MM-yyyy"
$ws.Cells.Item(20, 3).Value = "אחד נגד השני"
$ws.Cells.Item(20, 4).Value = 5
